$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("C2").Value = 10.5
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 9
$ws.Range("C5").Value = 16

# Add a new blank row 6 with the same style as row 5 (data rows)
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)
$ws.Range("A6:C6").ClearContents()

# Update selection to match target
$ws.Range("B3").Select()
